$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Charts to Make")

# Datasets for state landmass (row 8-9, "X-Axis title" = Adjusted Carbon Offset by Landmass)
# and state population (row 10-12) have been found, so mark the "Completed?2"
# column (G) as "Yes" for those rows, copying the same green-fill formatting
# already used on the completed rows above (G6:G7).
$ws.Range("G6").Copy() | Out-Null
$ws.Range("G8:G12").PasteSpecial(-4122) | Out-Null
$ws.Range("G8:G12").Value = "Yes"
$excel.CutCopyMode = 0

# Update the remembered selection on both sheets to reflect where the
# author was working when they saved.
$ws.Activate()
$ws.Range("F22").Select() | Out-Null

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Activate()
$wsMeta.Range("B17").Select() | Out-Null

$ws.Activate()
